# Generate Report for Handoff
#
# Two additional source files were discovered during the handoff run:
#   31257c33-f7a2-4fe7-a75c-cc20d410a7a7.md  (inserted between the existing
#                                              9e56c837... and c4c276b7... rows)
#   f7d9e043-725a-4f6d-826b-cc4a1e3ca2ef.md  (appended at the end)
#
# Each of the three tables (Overview, zh-cn, de-de) grows from 2 data rows
# to 4 data rows to track them.

$wb = $excel.ActiveWorkbook

$uuidNew1 = "31257c33-f7a2-4fe7-a75c-cc20d410a7a7"
$uuidOld  = "c4c276b7-2d3b-4581-9296-39aaf487959f"
$uuidNew2 = "f7d9e043-725a-4f6d-826b-cc4a1e3ca2ef"

$hashNew1 = "5b53b63bf5e96d352b4fea05a486da131b436d6a"
$hashOld  = "9270e297d57f2cc0ca1d3ad72a094f654c1b207d"
$hashNew2 = "1ab74bf8ced2d1188ae667fc6970bfa1a7c67209"

$srcCommitOld  = "c01c7336d69b14f40236c3385135ec581accc2c8"
$srcCommitNew1 = "a14f6ad6e1a2c7e0c5a2b3d4e5f6d41fbfe6b3c2"
$srcCommitNew2 = "2f3a4b5c6d7e8f90a1b2c3d4e5f6a77a8b9c0d1e"

$dtFmt = "yyyy-mm-dd HH:mm:ss"

# =======================================================================
# Sheet "Overview" -- columns: File Name | Path And Name | Extension |
#                               Publish URL | zh-cn | de-de | Latest HO Xliff Generate Date
# =======================================================================
$wsO = $wb.Worksheets.Item("Overview")
$loO = $wsO.ListObjects.Item(1)
$loO.ListRows.Add() | Out-Null
$loO.ListRows.Add() | Out-Null

function Set-OverviewRow($r, $uuid, $status, $dt) {
    $wsO.Range("A$r").Value = "$uuid.md"
    $wsO.Range("B$r").Value = "e2e\$uuid.md"
    $wsO.Range("C$r").Value = ".md"
    $wsO.Range("E$r").Value = $status
    $wsO.Range("F$r").Value = $status
    $wsO.Range("G$r").Value = $dt
    $wsO.Range("G$r").NumberFormat = $dtFmt
}

Set-OverviewRow 3 $uuidNew1 "Ready for handoff" "2016-08-30 08:49:24"
Set-OverviewRow 4 $uuidOld  "Ready for handoff" "2016-08-30 08:47:42"
Set-OverviewRow 5 $uuidNew2 "Ready for handoff" "2016-08-30 08:49:24"

$wsO.Range("B3").Hyperlinks.Delete()
$wsO.Hyperlinks.Add($wsO.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$srcCommitNew1/e2e/$uuidNew1.md", "", "", "e2e\$uuidNew1.md") | Out-Null
$wsO.Hyperlinks.Add($wsO.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$srcCommitOld/e2e/$uuidOld.md", "", "", "e2e\$uuidOld.md") | Out-Null
$wsO.Hyperlinks.Add($wsO.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$srcCommitNew2/e2e/$uuidNew2.md", "", "", "e2e\$uuidNew2.md") | Out-Null

# =======================================================================
# Sheets "zh-cn" / "de-de" -- columns: Source File Name | File Extension |
#   Status | Source Path | Priority | Content Duplicate | Latest Handoff
#   File | Latest Handoff Datetime | Latest Target File | Latest Handback
#   File | Latest Handback DateTime | Reference Tokens | To be localized |
#   Dependency From | Has metadata | Error Detail
# =======================================================================
function Set-LangRow($ws, $r, $uuid, $status, $hash, $lang, $handoffDt) {
    $ws.Range("A$r").Value = "$uuid.md"
    $ws.Range("B$r").Value = ".md"
    $ws.Range("C$r").Value = $status
    $ws.Range("D$r").Value = "e2e"
    $ws.Range("E$r").Value = "ht"
    $ws.Range("F$r").Value = "False"
    $ws.Range("G$r").Value = "$uuid.$hash.$lang.xlf"
    $ws.Range("H$r").Value = $handoffDt
    $ws.Range("H$r").NumberFormat = $dtFmt
    $ws.Range("K$r").Value = "0001-01-01 00:00:00"
    $ws.Range("K$r").NumberFormat = $dtFmt
    $ws.Range("M$r").Value = "True"
    $ws.Range("O$r").Value = "False"
}

# ---- zh-cn ----
$wsZ = $wb.Worksheets.Item("zh-cn")
$loZ = $wsZ.ListObjects.Item(1)
$loZ.ListRows.Add() | Out-Null
$loZ.ListRows.Add() | Out-Null

Set-LangRow $wsZ 3 $uuidNew1 "Ready for handoff" $hashNew1 "zh-cn" "2016-08-30 08:49:19"
Set-LangRow $wsZ 4 $uuidOld  "Ready for handoff" $hashOld  "zh-cn" "2016-08-30 08:47:37"
Set-LangRow $wsZ 5 $uuidNew2 "Ready for handoff" $hashNew2 "zh-cn" "2016-08-30 08:49:19"

$wsZ.Range("A3").Hyperlinks.Delete()
$wsZ.Hyperlinks.Add($wsZ.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$srcCommitNew1/e2e/$uuidNew1.md", "", "", "$uuidNew1.md") | Out-Null
$wsZ.Hyperlinks.Add($wsZ.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$srcCommitOld/e2e/$uuidOld.md", "", "", "$uuidOld.md") | Out-Null
$wsZ.Hyperlinks.Add($wsZ.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$srcCommitNew2/e2e/$uuidNew2.md", "", "", "$uuidNew2.md") | Out-Null

# ---- de-de ----
$wsD = $wb.Worksheets.Item("de-de")
$loD = $wsD.ListObjects.Item(1)
$loD.ListRows.Add() | Out-Null
$loD.ListRows.Add() | Out-Null

Set-LangRow $wsD 3 $uuidNew1 "Ready for handoff" $hashNew1 "de-de" "2016-08-30 08:49:24"
Set-LangRow $wsD 4 $uuidOld  "Ready for handoff" $hashOld  "de-de" "2016-08-30 08:47:42"
Set-LangRow $wsD 5 $uuidNew2 "Ready for handoff" $hashNew2 "de-de" "2016-08-30 08:49:24"

$wsD.Range("A3").Hyperlinks.Delete()
$wsD.Hyperlinks.Add($wsD.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$srcCommitNew1/e2e/$uuidNew1.md", "", "", "$uuidNew1.md") | Out-Null
$wsD.Hyperlinks.Add($wsD.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$srcCommitOld/e2e/$uuidOld.md", "", "", "$uuidOld.md") | Out-Null
$wsD.Hyperlinks.Add($wsD.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$srcCommitNew2/e2e/$uuidNew2.md", "", "", "$uuidNew2.md") | Out-Null

"Report updated: added $uuidNew1 and $uuidNew2 across Overview/zh-cn/de-de"
